$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing the existing data rows (2..16) down
# to (3..17) to make room for a new weekly price record.
$ws.Rows(2).Insert()
$ws.Range("A2:T2").ClearFormats()

# Columns A, B, C, E, F, G, H, I, J, K, Q, R, T hold the same constant values
# on every data row for this market/product, so copy them straight from the
# row below (which now carries what used to be row 2's data).
$ws.Cells.Item(2, 1).Value = $ws.Cells.Item(3, 1).Value2
$ws.Cells.Item(2, 2).Value = $ws.Cells.Item(3, 2).Value2
$ws.Cells.Item(2, 3).Value = $ws.Cells.Item(3, 3).Value2
$ws.Cells.Item(2, 5).Value = $ws.Cells.Item(3, 5).Value2
$ws.Cells.Item(2, 6).Value = $ws.Cells.Item(3, 6).Value2
$ws.Cells.Item(2, 7).Value = $ws.Cells.Item(3, 7).Value2
$ws.Cells.Item(2, 8).Value = $ws.Cells.Item(3, 8).Value2
$ws.Cells.Item(2, 9).Value = $ws.Cells.Item(3, 9).Value2
$ws.Cells.Item(2, 10).Value = $ws.Cells.Item(3, 10).Value2
$ws.Cells.Item(2, 11).Value = $ws.Cells.Item(3, 11).Value2
$ws.Cells.Item(2, 17).Value = $ws.Cells.Item(3, 17).Value2
$ws.Cells.Item(2, 18).Value = $ws.Cells.Item(3, 18).Value2
$ws.Cells.Item(2, 20).Value = $ws.Cells.Item(3, 20).Value2

# Date column keeps the same number format as the rest of column D.
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 4).Value = 44910

# New record-specific values for this week's entry.
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 200
$ws.Cells.Item(2, 14).Value = 7500
$ws.Cells.Item(2, 15).Value = 8000
$ws.Cells.Item(2, 16).Value = 7750
$ws.Cells.Item(2, 19).Value = 3875

$ws.Range("A1").Select()
